# BOM.xlsx update: "schematics finished (except gps)"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 12: canbus driver quantity 1 -> 3 (total recalculates via existing formula) ---
$ws.Range("D12").Value = 3

# --- Row 18: quartz (581-CX3225SB16D0FLJ) pricing filled in ---
$ws.Range("G18").Font.Color = 0
$ws.Range("G18").Value = 1.17
$ws.Range("H18").Font.Color = 0
$ws.Range("H18").Value = "Sylvestre van Kappel"
$ws.Range("I18").Formula = "=G18*D18"

# --- Row 19: common mode choke (871-B82793C104N201) new component line ---
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = "871-B82793C104N201 "
$ws.Range("F19").Font.Color = 0
$ws.Range("F19").Value = "common mode choke"
$ws.Range("G19").Font.Color = 0
$ws.Range("G19").Value = 2.0699999999999998
$ws.Range("H19").Font.Color = 0
$ws.Range("H19").Value = "Sylvestre van Kappel"
$ws.Range("I19").Formula = "=G19*D19"

# --- Row 20: self (495-TCK-141) new component line ---
$ws.Range("D20").Font.Color = 0
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = "495-TCK-141 "
$ws.Range("F20").Font.Color = 0
$ws.Range("F20").Value = "self"
$ws.Range("G20").Font.Color = 0
$ws.Range("G20").Value = 6.23
$ws.Range("H20").Font.Color = 0
$ws.Range("H20").Value = "Sylvestre van Kappel"

# --- Row 21: empty placeholder line, but "Qui" column now filled ---
$ws.Range("D21").Font.Color = 0
$ws.Range("H21").Font.Color = 0
$ws.Range("H21").VerticalAlignment = -4108
$ws.Range("H21").Value = "Sylvestre van Kappel"

# --- Row 22: "Qui" column filled ---
$ws.Range("H22").Font.Color = 0
$ws.Range("H22").VerticalAlignment = -4108
$ws.Range("H22").Value = "Sylvestre van Kappel"

# --- Row 23: blank spacer line, its old total cell is removed ---
$ws.Range("I23").ClearContents()

# --- Row 24: level shiffter reference swapped to 296-23759-6-ND, price 1.5, formula now points at D23 ---
$ws.Range("E24").Value = "296-23759-6-ND"
$ws.Range("G24").Value = 1.5
$ws.Range("H24").Font.Color = 0
$ws.Range("H24").VerticalAlignment = -4108
$ws.Range("I24").Formula = "=G24*D23"

# --- Row 25: becomes the canbus-driver-reference "level shiffter" line (qty 2), formula points at D24 ---
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = "296-21929-1-ND"
$ws.Range("F25").Value = "level shiffter"
$ws.Range("G25").Value = 1.03
$ws.Range("I25").Formula = "=G25*D24"

# --- Row 26: USB-UART converter line moves down here (qty 3), formula points at D25 ---
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = "1568-1504-ND"
$ws.Range("F26").Value = "USB-UART converter"
$ws.Range("G26").Value = 14.85
$ws.Range("H26").Value = "Sylvestre van Kappel"
$ws.Range("I26").Formula = "=G26*D25"

# --- Row 27: new diode de protection (497-13262-1-ND) line, formula points at D26 ---
$ws.Range("E27").Value = "497-13262-1-ND"
$ws.Range("F27").Value = "diode de protection"
$ws.Range("G27").Value = 0.41
$ws.Range("H27").Value = "Sylvestre van Kappel"
$ws.Range("I27").Formula = "=G27*D26"

# --- selection cursor now resting on D26 ---
$ws.Range("D26").Select()
